$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("AA4").Value = "Payable"
$ws.Range("AB4").Value = 5000
$ws.Range("AA17").Value = "Paid"
$ws.Range("AB17").Value = 500
$ws.Range("AA55").Value = "Paid"
$ws.Range("AB55").Value = 300
$ws.Range("AA73").Value = "Payable"
$ws.Range("AB73").Value = 5000
$ws.Range("AA102").Value = "Payable"
$ws.Range("AB102").Value = 200
$ws.Range("AA103").Value = "Payable"
$ws.Range("AB103").Value = 300
$ws.Range("AA126").Value = "Payable"
$ws.Range("AB126").Value = 300
$ws.Range("AA132").Value = "Payable"
$ws.Range("AB132").Value = 100
$ws.Range("AA134").Value = "Payable"
$ws.Range("AB134").Value = 300
$ws.Range("AA136").Value = "Payable"
$ws.Range("AB136").Value = 5000
$ws.Range("AA140").Value = "Payable"
$ws.Range("AB140").Value = 5000
$ws.Range("AA141").Value = "Payable"
$ws.Range("AB141").Value = 5000
$ws.Range("AA142").Value = "Payable"
$ws.Range("AB142").Value = 5000
$ws.Range("AA143").Value = "Payable"
$ws.Range("AB143").Value = 5000
$ws.Range("AA144").Value = "Payable"
$ws.Range("AB144").Value = 300
$ws.Range("AA145").Value = "Payable"
$ws.Range("AB145").Value = 10000
$ws.Range("AA146").Value = "Payable"
$ws.Range("AB146").Value = 10000
$ws.Range("AA147").Value = "Payable"
$ws.Range("AB147").Value = 5000
$ws.Range("AA148").Value = "Payable"
$ws.Range("AB148").Value = 5000
$ws.Range("AA150").Value = "Payable"
$ws.Range("AB150").Value = 10000
$ws.Range("AA153").Value = "Payable"
$ws.Range("AB153").Value = 300
$ws.Range("AA155").Value = "Payable"
$ws.Range("AB155").Value = 5000
$ws.Range("AA156").Value = "Payable"
$ws.Range("AB156").Value = 5000
$ws.Range("AA163").Value = "Payable"
$ws.Range("AB163").Value = 5000
$ws.Range("AA164").Value = "Payable"
$ws.Range("AB164").Value = 5000
$ws.Range("AA171").Value = "Payable"
$ws.Range("AB171").Value = 5000
$ws.Range("AA174").Value = "Payable"
$ws.Range("AB174").Value = 100
$ws.Range("AA176").Value = "Payable"
$ws.Range("AB176").Value = 100
$ws.Range("AA181").Value = "Payable"
$ws.Range("AB181").Value = 5000
$ws.Range("AA191").Value = "Payable"
$ws.Range("AB191").Value = 5000
$ws.Range("AA196").Value = "Payable"
$ws.Range("AB196").Value = 100
$ws.Range("AA201").Value = "Payable"
$ws.Range("AB201").Value = 300
$ws.Range("AA204").Value = "Payable"
$ws.Range("AB204").Value = 5000
$ws.Range("AA206").Value = "Payable"
$ws.Range("AB206").Value = 5000
$ws.Range("AA209").Value = "Paid"
$ws.Range("AB209").Value = 100
$ws.Range("AA211").Value = "Payable"
$ws.Range("AB211").Value = 5000
$ws.Range("AA214").Value = "Payable"
$ws.Range("AB214").Value = 5000
$ws.Range("AA215").Value = "Payable"
$ws.Range("AB215").Value = 5000
$ws.Range("AA220").Value = "Payable"
$ws.Range("AB220").Value = 5000
$ws.Range("AA221").Value = "Payable"
$ws.Range("AB221").Value = 300
$ws.Range("AA223").Value = "Payable"
$ws.Range("AB223").Value = 10000
$ws.Range("AA239").Value = "Payable"
$ws.Range("AB239").Value = 5000
$ws.Range("AA251").Value = "Payable"
$ws.Range("AB251").Value = 5000
$ws.Range("AA253").Value = "Paid"
$ws.Range("AB253").Value = 300
$ws.Range("AA255").Value = "Payable"
$ws.Range("AB255").Value = 5000
$ws.Range("AA257").Value = "Payable"
$ws.Range("AB257").Value = 5000
$ws.Range("AA258").Value = "Paid"
$ws.Range("AB258").Value = 50
$ws.Range("AA259").Value = "Paid"
$ws.Range("AB259").Value = 50
$ws.Range("AA260").Value = "Paid"
$ws.Range("AB260").Value = 2000
$ws.Range("AA261").Value = "Paid"
$ws.Range("AB261").Value = 500
$ws.Range("AA262").Value = "Paid"
$ws.Range("AB262").Value = 500
$ws.Range("AA263").Value = "Paid"
$ws.Range("AB263").Value = 200
$ws.Range("AA264").Value = "Paid"
$ws.Range("AB264").Value = 50
$ws.Range("AA265").Value = "Paid"
$ws.Range("AB265").Value = 100
$ws.Range("AA266").Value = "Paid"
$ws.Range("AB266").Value = 50
$ws.Range("AA267").Value = "Paid"
$ws.Range("AB267").Value = 50
$ws.Range("AA268").Value = "Paid"
$ws.Range("AB268").Value = 50
$ws.Range("AA269").Value = "Paid"
$ws.Range("AB269").Value = 2000
$ws.Range("AA270").Value = "Paid"
$ws.Range("AB270").Value = 500
$ws.Range("AA273").Value = "Payable"
$ws.Range("AB273").Value = 300
